$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "TRANSPORTE"
$ws.Range("C8").Value = "R$ 49.780"

$ws.Range("C8").Select()
